# Append two new booking rows (9 and 10) to the Bookings sheet.
# Columns: A Booking ID, B Date, C Customer Name, D Email, E Phone,
#          F Guests, G Plan, H Plan Price, I Total Price, J Status,
#          K Booking Date, L Special Requests
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 9
$ws.Range("A9").Value = "SNOW-355489"
$ws.Range("B9").NumberFormat = "@"
$ws.Range("B9").Value = "2/24/2026"
$ws.Range("C9").Value = "a"
$ws.Range("D9").Value = "a@b.com"
$ws.Range("E9").Value = "swee3"
$ws.Range("F9").Value = 2
$ws.Range("G9").Value = "Family Ski Package"
$ws.Range("H9").Value = 32000
$ws.Range("I9").Value = 64000
$ws.Range("J9").Value = "Confirmed"
$ws.Range("K9").NumberFormat = "@"
$ws.Range("K9").Value = "2/23/2026"
$ws.Range("L9").Value = ""

# Row 10
$ws.Range("A10").Value = "SNOW-355489"
$ws.Range("B10").NumberFormat = "@"
$ws.Range("B10").Value = "2/25/2026"
$ws.Range("C10").Value = "a"
$ws.Range("D10").Value = "a@b.com"
$ws.Range("E10").Value = "swee3"
$ws.Range("F10").Value = 2
$ws.Range("G10").Value = "Family Ski Package"
$ws.Range("H10").Value = 32000
$ws.Range("I10").Value = 64000
$ws.Range("J10").Value = "Confirmed"
$ws.Range("K10").NumberFormat = "@"
$ws.Range("K10").Value = "2/23/2026"
$ws.Range("L10").Value = ""
